# The workbook's sheets each hold a pasted statsmodels OLS summary in cell
# B2. The summary text embeds the timestamp of the run ("Date:" / "Time:"
# lines). This re-run moved from Thu 02 Jan 2020 20:48:37 to
# Sun 05 Jan 2020 21:22:15 - update that substring on every sheet, leaving
# the rest of each summary block untouched.

$wb = $excel.ActiveWorkbook

$oldDate = "Thu, 02 Jan 2020"
$newDate = "Sun, 05 Jan 2020"
$oldTime = "20:48:37"
$newTime = "21:22:15"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2
    if ($text -ne $null -and $text.GetType().Name -eq "String" -and $text.Contains("OLS Regression Results")) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        $cell.Value2 = $updated
    }
}
